$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = -2
$ws.Range("F9").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("F16").Value = -3
$ws.Range("F20").Value = -8
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = 0
$ws.Range("F31").Value = -1
$ws.Range("F33").Value = -1
$ws.Range("F35").Value = -3
$ws.Range("F42").Value = -6
$ws.Range("F45").Value = -5
$ws.Range("F46").Value = -3
$ws.Range("F47").Value = -5
$ws.Range("F51").Value = -2
$ws.Range("F72").Value = -3
$ws.Range("F73").Value = 6
$ws.Range("F76").Value = -2
$ws.Range("F78").Value = -5
